$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 83: new arrival entry (Monday, Jan 16)
$ws.Range("A83").Value = 82
$ws.Range("B83").Value = "Monday, Jan 16"
$ws.Range("C83").Value = "6:03 AM"
$ws.Range("D83").Value = "X7542"
$ws.Range("E83").Value = "New York"
$ws.Range("F83").Value = "(JFK)"
$ws.Range("G83").Value = "Challenge Airlines "
$ws.Range("H83").Value = "B744"
$ws.Range("I83").Value = "(OO-ACE)"
$ws.Range("J83").Value = "5:53 AM"
$ws.Range("L83").Value = "0 hours, -10 minutes"

# Row 84: new arrival entry (Monday, Jan 16)
$ws.Range("A84").Value = 83
$ws.Range("B84").Value = "Monday, Jan 16"
$ws.Range("C84").Value = "10:10 AM"
$ws.Range("D84").Value = "FR8224"
$ws.Range("E84").Value = "Bristol"
$ws.Range("F84").Value = "(BRS)"
$ws.Range("G84").Value = "Ryanair "
$ws.Range("H84").Value = "B738"
$ws.Range("I84").Value = "(EI-EKK)"
$ws.Range("J84").Value = "9:51 AM"
$ws.Range("L84").Value = "0 hours, -19 minutes"
